$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.005.52"
$ws.Range("E2").Value = '  +2.79%  '

$ws.Range("D3").Value = "'3.200.01"
$ws.Range("E3").Value = '  +1.52%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = "'536.81"
$ws.Range("E5").Value = '  +0.04%  '

$ws.Range("D6").Value = "'145.29"
$ws.Range("E6").Value = '  +4.19%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("E8").Value = '  +3.57%  '

$ws.Range("E9").Value = '  +0.67%  '

$ws.Range("E10").Value = '  +3.60%  '

$ws.Range("D11").Value = "'0.433"
$ws.Range("E11").Value = '  +3.11%  '

$ws.Range("D12").Value = "'3.755.81"
$ws.Range("E12").Value = '  +1.58%  '

$ws.Range("E13").Value = '  -0.93%  '

$ws.Range("D14").Value = "'26.07"
$ws.Range("E14").Value = '  +1.01%  '

$ws.Range("D15").Value = "'0.0000173"
$ws.Range("E15").Value = '  +2.65%  '

$ws.Range("D16").Value = "'60.103.92"

$ws.Range("D17").Value = "'3.192.17"
$ws.Range("E17").Value = '  +0.89%  '

$ws.Range("E18").Value = '  +0.80%  '

$ws.Range("D19").Value = "'13.19"
$ws.Range("E19").Value = '  +1.39%  '

$ws.Range("D20").Value = "'8.28"
$ws.Range("E20").Value = '  +0.53%  '

$ws.Range("D21").Value = "'378.88"
$ws.Range("E21").Value = '  +0.40%  '

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = '  -0.35%  '

$ws.Range("D23").Value = "'0.524"
$ws.Range("E23").Value = '  +2.04%  '

$ws.Range("D24").Value = "'70.08"
$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").Value = "'8.91"
$ws.Range("E25").Value = '  +10.23%  '

$ws.Range("D26").Value = "'0.169"
$ws.Range("E26").Value = '  +1.32%  '

$ws.Range("E27").Value = '  +1.90%  '

$ws.Range("D28").Value = "'0.0₃0898"
$ws.Range("E28").Value = '  +3.66%  '

$ws.Range("D29").Value = "'6.21"
$ws.Range("E29").Value = '  +1.17%  '

$ws.Range("D30").Value = "'1.90"
$ws.Range("E30").Value = '  +0.92%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = "'22.32"
$ws.Range("E31").Value = '  +2.11%  '

$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = "'5.42"
$ws.Range("E32").Value = '  +4.91%  '

$ws.Range("E33").Value = '  +3.58%  '

$ws.Range("D34").Value = "'6.66"
$ws.Range("E34").Value = '  +6.86%  '

$ws.Range("D35").Value = "'156.97"
$ws.Range("E35").Value = '  -2.46%  '

$ws.Range("E36").Value = '  -0.77%  '

$ws.Range("D37").Value = "'2.797.84"
$ws.Range("E37").Value = '  +5.84%  '

$ws.Range("D38").Value = "'25.52"
$ws.Range("E38").Value = '  +0.83%  '

$ws.Range("D39").Value = "'0.0704"
$ws.Range("E39").Value = '  +3.57%  '

$ws.Range("D40").Value = "'1.67"
$ws.Range("E40").Value = '  +0.89%  '

$ws.Range("D41").Value = "'4.24"
$ws.Range("E41").Value = '  +1.03%  '

$ws.Range("D42").Value = "'39.85"
$ws.Range("E42").Value = '  +3.12%  '

$ws.Range("D43").Value = "'0.0294"
$ws.Range("E43").Value = '  +5.10%  '

$ws.Range("E44").Value = '  +1.70%  '

$ws.Range("E45").Value = '  +3.27%  '

$ws.Range("D46").Value = "'3.245.03"
$ws.Range("E46").Value = '  +1.47%  '

$ws.Range("D47").Value = "'0.993"
$ws.Range("E47").Value = '  +1.61%  '

$ws.Range("D48").Value = "'6.17"
$ws.Range("E48").Value = '  -0.58%  '

$ws.Range("D49").Value = "'0.809"
$ws.Range("E49").Value = '  +7.03%  '

$ws.Range("D50").Value = "'20.62"
$ws.Range("E50").Value = '  +2.00%  '

$ws.Range("E51").Value = '  -0.02%  '
